$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 134
$ws.Range("I11").Value = 134
$ws.Range("K11").Value = 134
$ws.Range("M11").Value = 6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 11611.4
$ws.Range("J21").Value = 7019
$ws.Range("L21").Value = 7019
$ws.Range("N21").Value = -7955

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 11611.4
$ws.Range("J23").Value = 7019
$ws.Range("L23").Value = 7019
$ws.Range("N23").Value = -7487

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1682.0588
$ws.Range("J29").Value = 2315.8333
$ws.Range("L29").Value = 6947.499899999999
$ws.Range("N29").Value = -7509.499899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 428.83334
$ws.Range("I39").Value = 569.5
$ws.Range("J39").Value = 147.5
$ws.Range("K39").Value = 1708.5
$ws.Range("L39").Value = 442.5
$ws.Range("M39").Value = -1412.5
$ws.Range("N39").Value = -1034.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2300
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4480
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -5877
$ws.Range("N86").Value = -5046

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4480
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -29384
$ws.Range("N89").Value = -25232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3614.625
$ws.Range("I116").Value = 2343.6667
$ws.Range("K116").Value = 2343.6667
$ws.Range("M116").Value = 1098.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 83334344
$ws.Range("I61").Value = 100000820
$ws.Range("K61").Value = 100000820
$ws.Range("M61").Value = -100000608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2002.9166
$ws.Range("I122").Value = 1402.1
$ws.Range("K122").Value = 4206.299999999999
$ws.Range("M122").Value = -1756.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 83334344
$ws.Range("I136").Value = 100000820
$ws.Range("K136").Value = 300002460
$ws.Range("M136").Value = -299999910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 28318
$ws.Range("J139").Value = 28318
$ws.Range("L139").Value = 28318
$ws.Range("N139").Value = -38598

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 31132.5
$ws.Range("J141").Value = 31132.5
$ws.Range("L141").Value = 31132.5
$ws.Range("N141").Value = -41492.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 87648
$ws.Range("I22").Value = 147.33333
$ws.Range("K22").Value = 147.33333
$ws.Range("M22").Value = 202.66667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 24632.5
$ws.Range("J41").Value = 24632.5
$ws.Range("L41").Value = 24632.5
$ws.Range("N41").Value = -25488.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 26913.6
$ws.Range("J50").Value = 26913.6
$ws.Range("L50").Value = 26913.6
$ws.Range("N50").Value = -28163.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 23579.857
$ws.Range("J51").Value = 24176.5
$ws.Range("L51").Value = 24176.5
$ws.Range("N51").Value = -25648.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1525.6086
$ws.Range("I58").Value = 1198.2
$ws.Range("J58").Value = 2139.5
$ws.Range("K58").Value = 1198.2
$ws.Range("L58").Value = 2139.5
$ws.Range("M58").Value = -995.2
$ws.Range("N58").Value = -2545.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 23579.857
$ws.Range("J61").Value = 24176.5
$ws.Range("L61").Value = 24176.5
$ws.Range("N61").Value = -24872.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2611.8333
$ws.Range("I132").Value = 2332.5
$ws.Range("J132").Value = 2961
$ws.Range("K132").Value = 6997.5
$ws.Range("L132").Value = 8883
$ws.Range("M132").Value = -4467.5
$ws.Range("N132").Value = -13943

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 19232652
$ws.Range("I134").Value = 2043.1
$ws.Range("J134").Value = 83334690
$ws.Range("K134").Value = 6129.299999999999
$ws.Range("L134").Value = 250004070
$ws.Range("M134").Value = -3594.299999999999
$ws.Range("N134").Value = -250009140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1525.6086
$ws.Range("I136").Value = 1198.2
$ws.Range("J136").Value = 2139.5
$ws.Range("K136").Value = 3594.6
$ws.Range("L136").Value = 6418.5
$ws.Range("M136").Value = -1044.6
$ws.Range("N136").Value = -11518.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16396548
$ws.Range("I131").Value = 111111460
$ws.Range("J131").Value = 3582.923
$ws.Range("K131").Value = 333334380
$ws.Range("L131").Value = 10748.769
$ws.Range("M131").Value = -333329340
$ws.Range("N131").Value = -20828.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 5000
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5500
$ws.Range("J29").Value = 5500
$ws.Range("L29").Value = 5500
$ws.Range("N29").Value = -6080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 16999.75
$ws.Range("J46").Value = 21999.666
$ws.Range("L46").Value = 21999.666
$ws.Range("N46").Value = -22311.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 3093.75
$ws.Range("J54").Value = 3093.75
$ws.Range("L54").Value = 3093.75
$ws.Range("N54").Value = -3873.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3977.5334
$ws.Range("I132").Value = 4085.1
$ws.Range("J132").Value = 3762.4
$ws.Range("K132").Value = 12255.3
$ws.Range("L132").Value = 11287.2
$ws.Range("M132").Value = -9725.299999999999
$ws.Range("N132").Value = -16347.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1866.3334
$ws.Range("I7").Value = 1866.3334
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1866.3334
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1754.3334
$ws.Range("N7").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 713
$ws.Range("I22").Value = 439.4
$ws.Range("K22").Value = 439.4
$ws.Range("M22").Value = -144.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 713
$ws.Range("I27").Value = 439.4
$ws.Range("K27").Value = 439.4
$ws.Range("M27").Value = -332.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4402.636
$ws.Range("I46").Value = 489.8
$ws.Range("J46").Value = 7663.3335
$ws.Range("K46").Value = 489.8
$ws.Range("L46").Value = 7663.3335
$ws.Range("M46").Value = -301.8
$ws.Range("N46").Value = -8039.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1368
$ws.Range("I61").Value = 1399.5
$ws.Range("J61").Value = 1305
$ws.Range("K61").Value = 1399.5
$ws.Range("L61").Value = 1305
$ws.Range("M61").Value = -1197.5
$ws.Range("N61").Value = -1709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1368
$ws.Range("I113").Value = 1399.5
$ws.Range("J113").Value = 1305
$ws.Range("K113").Value = 1399.5
$ws.Range("L113").Value = 1305
$ws.Range("M113").Value = 770.5
$ws.Range("N113").Value = -5645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1866.3334
$ws.Range("I126").Value = 1866.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5599.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3129.0002
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1810.8
$ws.Range("I136").Value = 1358.2858
$ws.Range("J136").Value = 2866.6667
$ws.Range("K136").Value = 4074.8574
$ws.Range("L136").Value = 8600.000100000001
$ws.Range("M136").Value = -1524.8574
$ws.Range("N136").Value = -13700.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 34398
$ws.Range("J138").Value = 34398
$ws.Range("L138").Value = 34398
$ws.Range("N138").Value = -44678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 60509.8
$ws.Range("J140").Value = 60509.8
$ws.Range("L140").Value = 60509.8
$ws.Range("N140").Value = -70869.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6722.2666
$ws.Range("I132").Value = 9004
$ws.Range("K132").Value = 27012
$ws.Range("M132").Value = -24482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 999.2222
$ws.Range("I136").Value = 943.6667
$ws.Range("J136").Value = 1165.8889
$ws.Range("K136").Value = 2831.0001
$ws.Range("L136").Value = 3497.6667
$ws.Range("M136").Value = -281.0001000000002
$ws.Range("N136").Value = -8597.6667
